$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of rows 3 & 4 with rows 5 & 6 respectively
# (columns D, L, M, N, O, P, Q, R, S, T) while columns A, B, C, E-K stay the same.

$cols = @("D","L","M","N","O","P","Q","R","S","T")

foreach ($col in $cols) {
    $rangeA = "$col" + "3"
    $rangeB = "$col" + "5"
    $tmp = $ws.Range($rangeA).Value2
    $ws.Range($rangeA).Value2 = $ws.Range($rangeB).Value2
    $ws.Range($rangeB).Value2 = $tmp
}

foreach ($col in $cols) {
    $rangeA = "$col" + "4"
    $rangeB = "$col" + "6"
    $tmp = $ws.Range($rangeA).Value2
    $ws.Range($rangeA).Value2 = $ws.Range($rangeB).Value2
    $ws.Range($rangeB).Value2 = $tmp
}
